$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.791.51"
Set-TextValue $ws.Range("E2") "  -2.25%  "

Set-TextValue $ws.Range("D3") "3.633.43"
Set-TextValue $ws.Range("E3") "  +0.12%  "

Set-TextValue $ws.Range("E4") "  -0.08%  "

Set-TextValue $ws.Range("D5") "583.92"
Set-TextValue $ws.Range("E5") "  -2.27%  "

Set-TextValue $ws.Range("D6") "176.09"
Set-TextValue $ws.Range("E6") "  -3.43%  "

Set-TextValue $ws.Range("D7") "3.624.40"
Set-TextValue $ws.Range("E7") "  +0.14%  "

Set-TextValue $ws.Range("E8") "  +1.10%  "

Set-TextValue $ws.Range("E9") "  -0.11%  "

Set-TextValue $ws.Range("E10") "  -4.86%  "

Set-TextValue $ws.Range("D11") "6.85"
Set-TextValue $ws.Range("E11") "  +16.47%  "

Set-TextValue $ws.Range("D12") "0.608"
Set-TextValue $ws.Range("E12") "  -0.13%  "

Set-TextValue $ws.Range("E13") "  -3.83%  "

Set-TextValue $ws.Range("D14") "0.0000284"
Set-TextValue $ws.Range("E14") "  -2.07%  "

Set-TextValue $ws.Range("D15") "4.222.33"
Set-TextValue $ws.Range("E15") "  +0.18%  "

Set-TextValue $ws.Range("D16") "674.43"
Set-TextValue $ws.Range("E16") "  -5.01%  "

Set-TextValue $ws.Range("D17") "8.96"
Set-TextValue $ws.Range("E17") "  -0.25%  "

Set-TextValue $ws.Range("D18") "3.636.65"
Set-TextValue $ws.Range("E18") "  +0.93%  "

Set-TextValue $ws.Range("D19") "70.841.87"
Set-TextValue $ws.Range("E19") "  -2.37%  "

Set-TextValue $ws.Range("E20") "  -0.39%  "

Set-TextValue $ws.Range("E21") "  -4.26%  "

Set-TextValue $ws.Range("E22") "  -1.97%  "

Set-TextValue $ws.Range("E23") "  +0.92%  "

Set-TextValue $ws.Range("E24") "  -3.74%  "

Set-TextValue $ws.Range("D25") "100.08"
Set-TextValue $ws.Range("E25") "  -5.21%  "

Set-TextValue $ws.Range("E26") "  -2.83%  "

Set-TextValue $ws.Range("E27") "  -2.57%  "

Set-TextValue $ws.Range("E28") "  -0.10%  "

Set-TextValue $ws.Range("D29") "9.81"
Set-TextValue $ws.Range("E29") "  -2.45%  "

Set-TextValue $ws.Range("D30") "34.66"
Set-TextValue $ws.Range("E30") "  -2.45%  "

Set-TextValue $ws.Range("D31") "9.15"
Set-TextValue $ws.Range("E31") "  +0.21%  "

Set-TextValue $ws.Range("D32") "3.29"
Set-TextValue $ws.Range("E32") "  -5.91%  "

Set-TextValue $ws.Range("E33") "  +1.73%  "

Set-TextValue $ws.Range("E34") "  -6.39%  "

Set-TextValue $ws.Range("D35") "3.99"
Set-TextValue $ws.Range("E35") "  -4.77%  "

Set-TextValue $ws.Range("D36") "575.78"
Set-TextValue $ws.Range("E36") "  -3.25%  "

Set-TextValue $ws.Range("D37") "11.11"
Set-TextValue $ws.Range("E37") "  -2.26%  "

Set-TextValue $ws.Range("E38") "  -0.97%  "

Set-TextValue $ws.Range("D39") "58.45"
Set-TextValue $ws.Range("E39") "  -2.32%  "

Set-TextValue $ws.Range("D40") "1.00"
Set-TextValue $ws.Range("E40") "  +0.14%  "

Set-TextValue $ws.Range("D41") "3.558.22"
Set-TextValue $ws.Range("E41") "  -2.40%  "

Set-TextValue $ws.Range("D42") "0.0452"
Set-TextValue $ws.Range("E42") "  +0.37%  "

Set-TextValue $ws.Range("D43") "0.345"
Set-TextValue $ws.Range("E43") "  -1.61%  "

Set-TextValue $ws.Range("E44") "  -3.66%  "

Set-TextValue $ws.Range("D45") "34.32"
Set-TextValue $ws.Range("E45") "  -4.42%  "

Set-TextValue $ws.Range("D46") "0.0₃0732"
Set-TextValue $ws.Range("E46") "  -6.03%  "

Set-TextValue $ws.Range("D47") "2.68"
Set-TextValue $ws.Range("E47") "  -5.10%  "

Set-TextValue $ws.Range("D48") "2.94"
Set-TextValue $ws.Range("E48") "  +3.87%  "

Set-TextValue $ws.Range("E49") "  +0.96%  "

Set-TextValue $ws.Range("D50") "137.76"
Set-TextValue $ws.Range("E50") "  +3.06%  "

Set-TextValue $ws.Range("E51") "  -4.23%  "

